$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 4599.533
$ws.Range("I76").Value = 3010.3333
$ws.Range("K76").Value = 3010.3333
$ws.Range("M76").Value = -2695.3333
# Row 79
$ws.Range("H79").Value = 4599.533
$ws.Range("I79").Value = 3010.3333
$ws.Range("K79").Value = 3010.3333
$ws.Range("M79").Value = -1918.3333
# Row 82
$ws.Range("H82").Value = 2652.5715
$ws.Range("I82").Value = 1761.3334
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 5284.0002
$ws.Range("L82").Value = 24000
$ws.Range("M82").Value = -4878.0002
$ws.Range("N82").Value = -24812
# Row 85
$ws.Range("H85").Value = 2652.5715
$ws.Range("I85").Value = 1761.3334
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 5284.0002
$ws.Range("L85").Value = 24000
$ws.Range("M85").Value = -3880.0002
$ws.Range("N85").Value = -26808
# Row 98
$ws.Range("H98").Value = 1270
$ws.Range("I98").Value = 1281.6666
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 1281.6666
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = 216.3334
$ws.Range("N98").Value = -4196
# Row 113
$ws.Range("H113").Value = 59572.64
$ws.Range("I113").Value = 86193.47
$ws.Range("K113").Value = 86193.47
$ws.Range("M113").Value = -82939.47
# Row 116
$ws.Range("H116").Value = 2841.5
$ws.Range("I116").Value = 2117.5833
$ws.Range("K116").Value = 2117.5833
$ws.Range("M116").Value = 1324.4167
# Row 122
$ws.Range("H122").Value = 1270
$ws.Range("I122").Value = 1281.6666
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3844.9998
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1394.9998
$ws.Range("N122").Value = -8500
# Row 125
$ws.Range("H125").Value = 3666.6667
$ws.Range("I125").Value = 3666.6667
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 33000.0003
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -30540.0003
$ws.Range("N125").ClearContents()
# Row 135
$ws.Range("H135").Value = 1372.8889
$ws.Range("I135").Value = 1446
$ws.Range("J135").Value = 130
$ws.Range("K135").Value = 13014
$ws.Range("L135").Value = 1170
$ws.Range("M135").Value = -10479
$ws.Range("N135").Value = -6240
# Row 138
$ws.Range("H138").Value = 4631019
$ws.Range("I138").Value = 5953478.5
$ws.Range("J138").Value = 2412.1875
$ws.Range("K138").Value = 17860435.5
$ws.Range("L138").Value = 7236.5625
$ws.Range("M138").Value = -17855295.5
$ws.Range("N138").Value = -17516.5625
# Row 139
$ws.Range("H139").Value = 162500
$ws.Range("J139").Value = 162500
$ws.Range("L139").Value = 162500
$ws.Range("N139").Value = -172780
# Row 140
$ws.Range("H140").Value = 49726.668
$ws.Range("J140").Value = 49726.668
$ws.Range("L140").Value = 49726.668
$ws.Range("N140").Value = -60086.668

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8686.915999999999
$ws.Range("I32").Value = 8683.147000000001
$ws.Range("K32").Value = 8683.147000000001
$ws.Range("M32").Value = -8396.147000000001
# Row 97
$ws.Range("H97").Value = 4242.75
$ws.Range("I97").Value = 4125.1924
$ws.Range("J97").Value = 4752.1665
$ws.Range("K97").Value = 4125.1924
$ws.Range("L97").Value = 4752.1665
$ws.Range("M97").Value = -3629.1924
$ws.Range("N97").Value = -5744.1665
# Row 132
$ws.Range("H132").Value = 6252164.5
$ws.Range("I132").Value = 8930637
$ws.Range("J132").Value = 2393.5833
$ws.Range("K132").Value = 26791911
$ws.Range("L132").Value = 7180.749899999999
$ws.Range("M132").Value = -26789381
$ws.Range("N132").Value = -12240.7499

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 63
$ws.Range("H63").Value = 47647.332
$ws.Range("J63").Value = 47647.332
$ws.Range("L63").Value = 47647.332
$ws.Range("N63").Value = -49019.332
# Row 66
$ws.Range("H66").Value = 47647.332
$ws.Range("J66").Value = 47647.332
$ws.Range("L66").Value = 142941.996
$ws.Range("N66").Value = -149805.996
# Row 105
$ws.Range("H105").Value = 2497.2856
$ws.Range("I105").Value = 1407.2354
$ws.Range("J105").Value = 4181.909
$ws.Range("K105").Value = 1407.2354
$ws.Range("L105").Value = 4181.909
$ws.Range("M105").Value = 339.7646
$ws.Range("N105").Value = -7675.909
# Row 134
$ws.Range("H134").Value = 2314.7896
$ws.Range("I134").Value = 1477.3096
$ws.Range("J134").Value = 4659.7334
$ws.Range("K134").Value = 4431.9288
$ws.Range("L134").Value = 13979.2002
$ws.Range("M134").Value = -1896.9288
$ws.Range("N134").Value = -19049.2002

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 18
$ws.Range("H18").Value = 42684
$ws.Range("J18").Value = 42684
$ws.Range("L18").Value = 42684
$ws.Range("N18").Value = -43144
# Row 31
$ws.Range("H31").Value = 6174902
$ws.Range("I31").Value = 1762.4082
$ws.Range("J31").Value = 66671668
$ws.Range("K31").Value = 1762.4082
$ws.Range("L31").Value = 66671668
$ws.Range("M31").Value = -1467.4082
$ws.Range("N31").Value = -66672258
# Row 34
$ws.Range("H34").Value = 6174902
$ws.Range("I34").Value = 1762.4082
$ws.Range("J34").Value = 66671668
$ws.Range("K34").Value = 1762.4082
$ws.Range("L34").Value = 66671668
$ws.Range("M34").Value = -1560.4082
$ws.Range("N34").Value = -66672072
# Row 58
$ws.Range("H58").Value = 1336.2778
$ws.Range("I58").Value = 624.4138
$ws.Range("J58").Value = 4285.4287
$ws.Range("K58").Value = 624.4138
$ws.Range("L58").Value = 4285.4287
$ws.Range("M58").Value = -421.4138
$ws.Range("N58").Value = -4691.4287
# Row 122
$ws.Range("H122").Value = 1974.2188
$ws.Range("I122").Value = 1982
$ws.Range("K122").Value = 5946
$ws.Range("M122").Value = -3496
# Row 136
$ws.Range("H136").Value = 1336.2778
$ws.Range("I136").Value = 624.4138
$ws.Range("J136").Value = 4285.4287
$ws.Range("K136").Value = 1873.2414
$ws.Range("L136").Value = 12856.2861
$ws.Range("M136").Value = 676.7585999999999
$ws.Range("N136").Value = -17956.2861

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 690.2442
$ws.Range("I113").Value = 437.22858
$ws.Range("J113").Value = 863.8823
$ws.Range("K113").Value = 1311.68574
$ws.Range("L113").Value = 2591.6469
$ws.Range("M113").Value = 858.3142599999999
$ws.Range("N113").Value = -6931.6469
# Row 122
$ws.Range("H122").Value = 1817.1111
$ws.Range("J122").Value = 1040
$ws.Range("L122").Value = 9360
$ws.Range("N122").Value = -14260
# Row 125
$ws.Range("H125").Value = 3364.2144
$ws.Range("J125").Value = 3989.9
$ws.Range("L125").Value = 11969.7
$ws.Range("N125").Value = -21809.7
# Row 139
$ws.Range("H139").Value = 2864.3635
$ws.Range("I139").Value = 1303.125
$ws.Range("K139").Value = 3909.375
$ws.Range("M139").Value = 1230.625

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2373.9534
$ws.Range("I102").Value = 2692.2122
$ws.Range("J102").Value = 1323.7
$ws.Range("K102").Value = 2692.2122
$ws.Range("L102").Value = 1323.7
$ws.Range("M102").Value = -1070.2122
$ws.Range("N102").Value = -4567.7
# Row 126
$ws.Range("H126").Value = 4678.1177
$ws.Range("I126").Value = 2650
$ws.Range("J126").Value = 5784.364
$ws.Range("K126").Value = 7950
$ws.Range("L126").Value = 17353.092
$ws.Range("M126").Value = -5480
$ws.Range("N126").Value = -22293.092
# Row 132
$ws.Range("H132").Value = 3686.276
$ws.Range("I132").Value = 3053.5312
$ws.Range("J132").Value = 4465.0386
$ws.Range("K132").Value = 9160.5936
$ws.Range("L132").Value = 13395.1158
$ws.Range("M132").Value = -6630.5936
$ws.Range("N132").Value = -18455.1158
# Row 138
$ws.Range("H138").Value = 52800
$ws.Range("J138").Value = 52800
$ws.Range("L138").Value = 52800
$ws.Range("N138").Value = -63080

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5752.263
$ws.Range("I7").Value = 6618.5
$ws.Range("J7").Value = 5122.273
$ws.Range("K7").Value = 6618.5
$ws.Range("L7").Value = 5122.273
$ws.Range("M7").Value = -6506.5
$ws.Range("N7").Value = -5346.273
# Row 82
$ws.Range("H82").Value = 1709.875
$ws.Range("I82").Value = 1550.6666
$ws.Range("J82").Value = 2187.5
$ws.Range("K82").Value = 1550.6666
$ws.Range("L82").Value = 2187.5
$ws.Range("M82").Value = -1189.6666
$ws.Range("N82").Value = -2909.5
# Row 85
$ws.Range("H85").Value = 1709.875
$ws.Range("I85").Value = 1550.6666
$ws.Range("J85").Value = 2187.5
$ws.Range("K85").Value = 1550.6666
$ws.Range("L85").Value = 2187.5
$ws.Range("M85").Value = -302.6666
$ws.Range("N85").Value = -4683.5
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 126
$ws.Range("H126").Value = 5752.263
$ws.Range("I126").Value = 6618.5
$ws.Range("J126").Value = 5122.273
$ws.Range("K126").Value = 19855.5
$ws.Range("L126").Value = 15366.819
$ws.Range("M126").Value = -17385.5
$ws.Range("N126").Value = -20306.819

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 1250400
$ws.Range("I3").Value = 1250400
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1250400
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1250286
$ws.Range("N3").ClearContents()
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 126
$ws.Range("H126").Value = 3180.6365
$ws.Range("I126").Value = 2342.4443
$ws.Range("K126").Value = 7027.3329
$ws.Range("M126").Value = -4557.3329
# Row 136
$ws.Range("H136").Value = 719.2545
$ws.Range("I136").Value = 555.28
$ws.Range("J136").Value = 2359
$ws.Range("K136").Value = 1665.84
$ws.Range("L136").Value = 7077
$ws.Range("M136").Value = 884.1600000000001
$ws.Range("N136").Value = -12177

Write-Host "Applied all profit updates"